# Daily attendance processing - 2026-01-31 12:00:52
# Swap the order of names in column G for rows where the value is
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G is the 7th column
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
